$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")
$ws.Activate()

# --- Question 9 (rows 98:100): INDEX/MATCH versions of the Question 8 lookups ---

# FY17 row (98) uses columns E (pct) / F (rank)
$ws.Range("B98").Formula = '=INDEX($A$2:$A$52, MATCH($B$96, F$2:F$52,0))'
$ws.Range("C98").Formula = '=INDEX($E$2:$E$52, MATCH($B$96, F$2:F$52,0))'
$ws.Range("D98").Formula = '=INDEX($A$2:$A$52, MATCH($D$96, F$2:F$52,0))'
$ws.Range("E98").Formula = '=INDEX($E$2:$E$52, MATCH(D$96, F$2:F$52,0))'
$ws.Range("F98").Formula = '=INDEX($A$2:$A$52, MATCH($F$96, F$2:F$52,0))'
$ws.Range("G98").Formula = '=INDEX($E$2:$E$52, MATCH(F$96, F$2:F$52,0))'

# FY18 row (99) uses columns J (pct) / K (rank)
$ws.Range("B99").Formula = '=INDEX($A$2:$A$52, MATCH($B$96, K$2:K$52,0))'
$ws.Range("C99").Formula = '=INDEX($J$2:$J$52, MATCH($B$96, K$2:K$52,0))'
$ws.Range("D99").Formula = '=INDEX($A$2:$A$52, MATCH($D$96, K$2:K$52,0))'
$ws.Range("E99").Formula = '=INDEX($J$2:$J$52, MATCH($D$96, K$2:K$52,0))'
$ws.Range("F99").Formula = '=INDEX($A$2:$A$52, MATCH($F$96, K$2:K$52,0))'
$ws.Range("G99").Formula = '=INDEX($J$2:$J$52, MATCH($F$96, K$2:K$52,0))'

# FY19 row (100) uses columns O (pct) / P (rank)
$ws.Range("B100").Formula = '=INDEX($A$2:$A$52, MATCH($B$96, P$2:P$52,0))'
$ws.Range("C100").Formula = '=INDEX($O$2:$O$52, MATCH($B$96, P$2:P$52,0))'
$ws.Range("D100").Formula = '=INDEX($A$2:$A$52, MATCH($D$96, P$2:P$52,0))'
$ws.Range("E100").Formula = '=INDEX($O$2:$O$52, MATCH($D$96, P$2:P$52,0))'
$ws.Range("F100").Formula = '=INDEX($A$2:$A$52, MATCH($F$96, P$2:P$52,0))'
$ws.Range("G100").Formula = '=INDEX($O$2:$O$52, MATCH($F$96, P$2:P$52,0))'

# Match the percent-style formatting (0.00%) used by the equivalent Question 8 cells.
# (Set one cell at a time - applying NumberFormat to a multi-area Range only affects
# the first area in this environment.)
foreach ($addr in @("C98","E98","G98","C99","E99","G99","C100","E100","G100")) {
    $ws.Range($addr).NumberFormat = "0.00%"
}

# --- Window/selection state: user navigated to and selected D98 ---
$ws.Range("D98").Select()
